# "some lesson 2 preps" - add lesson 2 (multiples of 3/4, division) strings
# to the "en" localization sheet, ahead of the existing lesson 1 rows, and
# refresh a couple of lesson 1 strings that changed wording at the same time.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

$rows = @(
    @(22, 'mult2_title', 'Multiples of 2'),
    @(23, 'mult3_title', 'Multiples of 3'),
    @(24, 'mult4_title', 'Multiples of 4'),
    @(25, 'mult5_title', 'Multiples of 5'),
    @(26, 'mult6_title', 'Multiples of 6'),
    @(27, 'mult7_title', 'Multiples of 7'),
    @(28, 'mult8_title', 'Multiples of 8'),
    @(29, 'mult9_title', 'Multiples of 9'),
    @(30, 'mult10_title', 'Multiples of 10'),
    @(31, 'division', 'Division'),
    @(32, 'commutative_title', 'Commutative Property'),
    @(33, 'not_commutative', 'Not Commutative!'),
    @(34, 'lesson_1_intro_1', 'Before we begin our mission, let''s first learn about some tricks with multiplication!'),
    @(35, 'lesson_1_mult2_1', 'In multiples of two, the trick is to simply double the number.'),
    @(36, 'lesson_1_mult2_2', 'For example: 2 x 6 can be 6 + 6, which equals to 12.'),
    @(37, 'lesson_1_commutative_1', 'The commutative property means that multiplying numbers in any order gives the same answer.'),
    @(38, 'lesson_1_commutative_2', 'For example: 2 x 3, and 3 x 2, equal 6.'),
    @(39, 'lesson_1_commutative_3', 'With this trick, you only have to remember half the multiplication table!'),
    @(40, 'lesson_1_tutorial_1', 'Now banish these blobs by connecting them in the correct order using multiplication.'),
    @(41, 'lesson_1_tutorial_end_1', 'Excellent! You are now ready for the mission!'),
    @(42, 'lesson_2_mult3_1', 'In multiples of three: double the number, and then add the original number.'),
    @(43, 'lesson_2_mult3_2', 'For example: 3 x 6 can be expressed as (6 x 2) + 6, which becomes 12 + 6, giving you 18.'),
    @(44, 'lesson_2_mult4_1', 'In multiples of four: double the number, and then double it again.'),
    @(45, 'lesson_2_mult4_2', 'For example: 4 x 6, double 6 to get 12, and then double 12 to get 24.'),
    @(46, 'lesson_2_div_1', 'When it comes to division, think of it as the opposite of multiplication.'),
    @(47, 'lesson_2_div_2', 'Rearranging the equation, and replacing division with multiplication can help.'),
    @(48, 'lesson_2_div_3', 'For example: 21 ÷ 3 = ? can be rearranged to ? x 3 = 21. From here, we can deduce that 7 x 3 = 21.'),
    @(49, 'lesson_2_div_4', 'Unlike multiplication, division is not commutative. So the order of the numbers cannot be changed.'),
    @(50, 'lesson_2_tutorial_1', 'Now banish these blobs by connecting them in the correct order using division.'),
    @(51, 'lesson_2_tutorial_end_1', 'Excellent! You are now ready for the mission!')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

$ws.Range("B33").Select() | Out-Null